$d = $word.ActiveDocument

# --- 1. Insert "Master " before "Username" and before "Password" in the
#        existing paragraph that describes the file format. ---
$d.Content.Find.Execute("Username for that Account", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Master Username for that Account", 2) | Out-Null

$d.Content.Find.Execute("a Password, and a filename", $true, $false, $false, $false, $false, `
    $true, 1, $false, "a Master Password, and a filename", 2) | Out-Null

# --- 2. Locate that paragraph's index (the one that now contains
#        "Master Password") so we can append new paragraphs right after it. ---
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Master Password*") {
        $targetIndex = $i
    }
}

# --- 3. Add a blank paragraph, a "File Header: " paragraph, and a
#        paragraph with a tab followed by "Account Name, Master Password",
#        all inheriting the same indented formatting as the paragraph
#        they follow. ---
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphAfter()

$blankIndex = $targetIndex + 1
$d.Paragraphs.Item($blankIndex).Range.InsertParagraphAfter()

$headerIndex = $blankIndex + 1
$d.Paragraphs.Item($headerIndex).Range.Text = "File Header: "
$d.Paragraphs.Item($headerIndex).Range.InsertParagraphAfter()

$detailIndex = $headerIndex + 1
$d.Paragraphs.Item($detailIndex).Range.Text = "`tAccount Name, Master Password"
